# Convert all v1 (3-sheet: library_content / requirements / answers) workbooks
# to v2 (5-sheet: *_meta / *_content split) layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) library_content -> library_meta  (re-key rows, drop framework_*/tab rows)
# ---------------------------------------------------------------------------
$wsLib = $wb.Worksheets.Item(1)
$wsLib.Name = "library_meta"

# Clear the old content (rows 1-15, cols A-C) before rewriting.
$wsLib.Range("A1:C15").ClearContents()

$wsLib.Range("A1").Value = "type"
$wsLib.Range("B1").Value = "library"

$wsLib.Range("A2").Value = "urn"
$wsLib.Range("B2").Value = "urn:intuitem:risk:library:adobe-ccf-v5"

$wsLib.Range("A3").Value = "version"
$wsLib.Range("B3").NumberFormat = "@"
$wsLib.Range("B3").Value = "1"
$wsLib.Range("B3").ClearFormats()

$wsLib.Range("A4").Value = "locale"
$wsLib.Range("B4").Value = "en"

$wsLib.Range("A5").Value = "ref_id"
$wsLib.Range("B5").Value = "adobe-ccf-v5"

$wsLib.Range("A6").Value = "name"
$wsLib.Range("B6").Value = "Adobe CCF v5"

$wsLib.Range("A7").Value = "description"
$wsLib.Range("B7").Value = "Adobe Common Controls Framework (CCF) version 5`nhttps://www.adobe.com/trust/compliance/adobe-ccf.html"

$wsLib.Range("A8").Value = "copyright"
$wsLib.Range("B8").Value = "Creative Commons"

$wsLib.Range("A9").Value = "provider"
$wsLib.Range("B9").Value = "Adobe"

$wsLib.Range("A10").Value = "packager"
$wsLib.Range("B10").Value = "intuitem"

# ---------------------------------------------------------------------------
# 2) requirements -> requirements_meta (new small framework-metadata sheet)
#    + requirements_content (verbatim copy of the old big requirements table,
#    with the blank category-header cells trimmed)
# ---------------------------------------------------------------------------
$wsReq = $wb.Worksheets.Item(2)

# Duplicate the sheet first so the full requirements table survives under a
# new tab, placed immediately after the original.
$wsReq.Copy($null, $wsReq)
$wsReqContent = $wb.Worksheets.Item(3)
$wsReqContent.Name = "requirements_content"

# Trim the now-empty A/C/E/F cells on the 25 category-header rows (rows that
# only carry a depth=1 marker in B and the category name in D).
$reqLastRow = $wsReqContent.UsedRange.Rows.Count
for ($r = 2; $r -le $reqLastRow; $r++) {
    $bVal = $wsReqContent.Cells.Item($r, 2).Text
    $cVal = $wsReqContent.Cells.Item($r, 3).Text
    if ($bVal -eq "1" -and $cVal -eq "") {
        $wsReqContent.Cells.Item($r, 1).ClearContents()
        $wsReqContent.Cells.Item($r, 3).ClearContents()
        $wsReqContent.Cells.Item($r, 5).ClearContents()
        $wsReqContent.Cells.Item($r, 6).ClearContents()
    }
}

# Now repurpose the original "requirements" tab as the compact framework
# metadata sheet.
$wsReq.Name = "requirements_meta"
$wsReq.Range("A1:I343").ClearContents()

$wsReq.Range("A1").Value = "type"
$wsReq.Range("B1").Value = "framework"

$wsReq.Range("A2").Value = "base_urn"
$wsReq.Range("B2").Value = "urn:intuitem:risk:req_node:adobe-ccf-v5"

$wsReq.Range("A3").Value = "urn"
$wsReq.Range("B3").Value = "urn:intuitem:risk:framework:adobe-ccf-v5"

$wsReq.Range("A4").Value = "ref_id"
$wsReq.Range("B4").Value = "adobe-ccf-v5"

$wsReq.Range("A5").Value = "name"
$wsReq.Range("B5").Value = "Adobe CCF v5"

$wsReq.Range("A6").Value = "description"
$wsReq.Range("B6").Value = "Adobe Common Controls Framework (CCF) version 5`nhttps://www.adobe.com/trust/compliance/adobe-ccf.html"

$wsReq.Range("A7").Value = "answers_definition"
$wsReq.Range("B7").Value = "answers"

# ---------------------------------------------------------------------------
# 3) answers -> answers_meta (new small sheet) + answers_content (verbatim
#    copy of the old answers table)
# ---------------------------------------------------------------------------
$wsAns = $wb.Worksheets.Item(4)

$wsAns.Copy($null, $wsAns)
$wsAnsContent = $wb.Worksheets.Item(5)
$wsAnsContent.Name = "answers_content"

$wsAns.Name = "answers_meta"
$wsAns.Range("A1:C2").ClearContents()

$wsAns.Range("A1").Value = "type"
$wsAns.Range("B1").Value = "answers"

$wsAns.Range("A2").Value = "name"
$wsAns.Range("B2").Value = "answers"
